$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row at the top of the episode list (row 2), shifting the
# --- existing episodes down by one. This makes room for the "currently
# --- playing" episode as the new first entry.
$ws.Range("A2").EntireRow.Insert()

# Row-insert in Excel carries formatting down from the row above for the
# newly inserted cells; the diff shows the new row's B:E cells with no
# explicit style (just like every other data row), so strip that back off.
$ws.Range("B2:E2").ClearFormats()

# Column A keeps its bold/bordered "index" style (style index 1 in the
# original sheet) on every data row, including the new one - copy that
# formatting down from the row below (now row 3, which held the old row 2).
$ws.Range("A3").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Fill in the new "currently playing" episode in row 2.
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "Mon, 06 Dec 2021 16:30:17 GMT"
$ws.Range("C2").Value = "The World Ahead: Year three"
$ws.Range("D2").Value = "00:26:54"
$ws.Range("E2").Value = "https://sphinx.acast.com/theeconomistallaudio/theworldahead/theworldahead-yearthree/media.mp3"

# --- Renumber the index column for the rows that shifted down.
$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 2
$ws.Range("A5").Value = 3
$ws.Range("A6").Value = 4
$ws.Range("A7").Value = 5

# --- Append two more episodes at the bottom of the list (rows 8 and 9),
# --- matching the index-column style used by the rest of the data rows.
$ws.Range("A7").Copy()
$ws.Range("A8:A9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "Thu, 02 Dec 2021 11:05:06 GMT"
$ws.Range("C8").Value = "Roe blow? SCOTUS weighs abortion rights"
$ws.Range("D8").Value = "00:24:02"
$ws.Range("E8").Value = "https://sphinx.acast.com/theeconomistallaudio/theintelligencepodcast/roeblow-scotusweighsabortionrights/media.mp3"

$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "Wed, 01 Dec 2021 17:54:16 GMT"
$ws.Range("C9").Value = "Money Talks: Omicronomics"
$ws.Range("D9").Value = "00:31:01"
$ws.Range("E9").Value = "https://sphinx.acast.com/theeconomistallaudio/theeconomistmoneytalks/moneytalks-omicronomics/media.mp3"
